$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data (row 63),
# pushing all subsequent rows (previously 63-79) down by one (now 64-80).
$ws.Rows(63).Insert()

$ws.Range("A63").Value = 7
$ws.Range("B63").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C63").Value = "Ñuble"
$ws.Range("D63").Value = 44637
$ws.Range("E63").Value = 16
$ws.Range("F63").Value = 100112031
$ws.Range("G63").Value = "Poroto verde"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 120
$ws.Range("K63").Value = 29000
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = 29500
$ws.Range("N63").Value = "`$/saco 25 kilos"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 1180
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
